$d = $word.ActiveDocument
$paras = $d.Paragraphs

# --- Table 1 (תורה בקיאות): comment paragraph #5, score paragraph #8 ---
$comment1 = @"
במחיצת זאת למדנו על חומש "ויקרא", למדנו את ההלכות והאיסורים לעומק,
הייתה אוירת לימוד מצוינת.
חיה את תלמידה מקסימה, הרבה בהצלחה!
"@
$paras.Item(5).Range.Text = $comment1
$paras.Item(8).Range.Text = "85"

# --- Table 2 (היסטוריה): comment paragraph #15, score paragraph #18 ---
$comment2 = @"
במחצית זאת למדנו על היסטורית השואה, לכל אחת היתה משימה לעשות פרוייקט על השואה, כך שחפרנו עמוק בשורשים.
חיה את ילדה נפלאה, הרבה הצלחה!
"@
$paras.Item(15).Range.Text = $comment2
$paras.Item(18).Range.Text = "84"

# --- Table 3 (מתמטיקה): comment paragraph #25, score paragraph #28 ---
$comment3 = @"
במחצית זאת למדנו על תורת המיספרים הגדולה, התמקדנו על שברים, על תורת המעגל, רדיוס וקטרים, הרחבנו בנושא המשוואות בנעלם אחד,
שיננו לעצמינו את הכללים החשובים שנזכור לעתיד.
חיה את ילדה מצוינת, את מעולה שיהיה לך הרבה הצלחה להמשך!
"@
$paras.Item(25).Range.Text = $comment3
$paras.Item(28).Range.Text = "80"

# --- Table 4 (אנגלית): comment paragraph #35 (ends with a trailing blank line), score paragraph #38 ---
$comment4 = @"
במחצית זאת למדנו את נושא השיכות, והתמקדנו על הדקדוק והזמנים, הווה מושלם והווה פשוט,היתה אוירה טובה!
חיה את מצוינת, בהצלחה!

"@
$paras.Item(35).Range.Text = $comment4
$paras.Item(38).Range.Text = "87"
